$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column O entirely (its values were duplicated/obsolete data); this
# shifts the used range from A1:O25 down to A1:N25 and updates all row spans
# automatically.
$ws.Columns("O").Delete()

# Rows 17-25 also got their M/N figures revised (M: 290 -> 0, N: 100 -> 50).
for ($r = 17; $r -le 25; $r++) {
    $ws.Cells.Item($r, 13).Value = 0   # column M
    $ws.Cells.Item($r, 14).Value = 50  # column N
}

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("Q23").Select()
